$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 09:22:27"
$ws.Range("A3").Value = "Total filas: 146"
$ws.Range("A65").Value = "07:12:47"
$ws.Range("C65").Value = "27_EL RETIRO"
$ws.Range("D65").Value = 24
$ws.Range("A66").Value = "06:46:37"
$ws.Range("C66").Value = "17X38_ROMERO"
$ws.Range("D66").Value = 50
$ws.Range("A75").Value = "06:53:56"
$ws.Range("C75").Value = "23_HERNANDEZ"
$ws.Range("D75").Value = 66
$ws.Range("A76").Value = "07:50:33"
$ws.Range("C76").Value = "11_ETCHEVERRY"
$ws.Range("D76").Value = 9
$ws.Range("C91").Value = "10_OLMOS"
$ws.Range("C92").Value = "215A_EL PATO"
$ws.Range("A103").Value = "08:52:13"
$ws.Range("C103").Value = "17X38_ROMERO"
$ws.Range("D103").Value = 11
$ws.Range("A104").Value = "08:29:58"
$ws.Range("C104").Value = "23_HERNANDEZ"
$ws.Range("D104").Value = 34
$ws.Range("A110").Value = "08:52:13"
$ws.Range("C110").Value = "11_ETCHEVERRY"
$ws.Range("D110").Value = 23
$ws.Range("A111").Value = "07:38:30"
$ws.Range("C111").Value = "27_EL RETIRO"
$ws.Range("D111").Value = 97
$ws.Range("A116").Value = "09:22:27"
$ws.Range("C116").Value = "10_OLMOS"
$ws.Range("D116").Value = 4
$ws.Range("A117").Value = "08:29:58"
$ws.Range("B117").Value = "09:26"
$ws.Range("C117").Value = "23_HERNANDEZ"
$ws.Range("D117").Value = 57
$ws.Range("A118").Value = "08:40:53"
$ws.Range("B118").Value = "09:28"
$ws.Range("D118").Value = 48
$ws.Range("A119").Value = "08:52:13"
$ws.Range("C119").Value = "10_OLMOS"
$ws.Range("D119").Value = 37
$ws.Range("A120").Value = "08:40:53"
$ws.Range("B120").Value = "09:29"
$ws.Range("D120").Value = 49
$ws.Range("A122").Value = "09:22:27"
$ws.Range("B122").Value = "09:33"
$ws.Range("C122").Value = "23_HERNANDEZ"
$ws.Range("D122").Value = 11
$ws.Range("A123").Value = "09:22:27"
$ws.Range("B123").Value = "09:34"
$ws.Range("C123").Value = "15_ABASTO"
$ws.Range("D123").Value = 12
$ws.Range("A124").Value = "09:22:27"
$ws.Range("B124").Value = "09:41"
$ws.Range("C124").Value = "16_SANTA ANA"
$ws.Range("D124").Value = 19
$ws.Range("A125").Value = "08:40:53"
$ws.Range("B125").Value = "09:44"
$ws.Range("D125").Value = 64
$ws.Range("A126").Value = "08:52:13"
$ws.Range("B126").Value = "09:45"
$ws.Range("C126").Value = "14_ABASTO"
$ws.Range("D126").Value = 53
$ws.Range("A127").Value = "08:29:58"
$ws.Range("B127").Value = "09:48"
$ws.Range("D127").Value = 79
$ws.Range("A128").Value = "08:10:38"
$ws.Range("B128").Value = "09:49"
$ws.Range("C128").Value = "15_ABASTO"
$ws.Range("D128").Value = 99
$ws.Range("A129").Value = "08:29:58"
$ws.Range("B129").Value = "09:50"
$ws.Range("D129").Value = 81
$ws.Range("A130").Value = "09:22:27"
$ws.Range("B130").Value = "09:51"
$ws.Range("C130").Value = "16_P MOR-SANTA ANA"
$ws.Range("D130").Value = 29
$ws.Range("A131").Value = "09:22:27"
$ws.Range("B131").Value = "09:56"
$ws.Range("C131").Value = "10_OLMOS"
$ws.Range("D131").Value = 34
$ws.Range("A132").Value = "08:40:53"
$ws.Range("B132").Value = "10:03"
$ws.Range("D132").Value = 83
$ws.Range("A133").Value = "09:22:27"
$ws.Range("B133").Value = "10:03"
$ws.Range("C133").Value = "23_HERNANDEZ"
$ws.Range("D133").Value = 41
$ws.Range("A134").Value = "09:22:27"
$ws.Range("B134").Value = "10:04"
$ws.Range("C134").Value = "215C_EL PATO"
$ws.Range("D134").Value = 42
$ws.Range("A135").Value = "09:22:27"
$ws.Range("B135").Value = "10:08"
$ws.Range("C135").Value = "11_ETCHEVERRY"
$ws.Range("D135").Value = 46
$ws.Range("B136").Value = "10:09"
$ws.Range("C136").Value = "11_ETCHEVERRY"
$ws.Range("D136").Value = 77
$ws.Range("A137").Value = "08:40:53"
$ws.Range("B137").Value = "10:18"
$ws.Range("C137").Value = "17_ROMERO"
$ws.Range("D137").Value = 98
$ws.Range("A138").Value = "09:22:27"
$ws.Range("B138").Value = "10:19"
$ws.Range("C138").Value = "17_ROMERO"
$ws.Range("D138").Value = 57
$ws.Range("A139").Value = "09:22:27"
$ws.Range("B139").Value = "10:20"
$ws.Range("C139").Value = "10_OLMOS"
$ws.Range("D139").Value = 58
$ws.Range("A140").Value = "09:22:27"
$ws.Range("B140").Value = "10:32"
$ws.Range("C140").Value = "14_ABASTO"
$ws.Range("D140").Value = 70
$ws.Range("B141").Value = "10:33"
$ws.Range("C141").Value = "14_ABASTO"
$ws.Range("D141").Value = 101
$ws.Range("A142").Value = "09:22:27"
$ws.Range("B142").Value = "10:34"
$ws.Range("C142").Value = "15_ABASTO"
$ws.Range("D142").Value = 72
$ws.Range("E142").Value = "LP1912"
$ws.Range("A143").Value = "09:22:27"
$ws.Range("B143").Value = "10:40"
$ws.Range("C143").Value = "16_SANTA ANA"
$ws.Range("D143").Value = 78
$ws.Range("E143").Value = "LP1912"
$ws.Range("A144").Value = "09:22:27"
$ws.Range("B144").Value = "10:44"
$ws.Range("C144").Value = "10_OLMOS"
$ws.Range("D144").Value = 82
$ws.Range("E144").Value = "LP1912"
$ws.Range("A145").Value = "09:22:27"
$ws.Range("B145").Value = "10:51"
$ws.Range("C145").Value = "16_P MOR-SANTA ANA"
$ws.Range("D145").Value = 89
$ws.Range("E145").Value = "LP1912"
$ws.Range("A146").Value = "09:22:27"
$ws.Range("B146").Value = "10:56"
$ws.Range("C146").Value = "27_EL RETIRO"
$ws.Range("D146").Value = 94
$ws.Range("E146").Value = "LP1912"
$ws.Range("A147").Value = "09:22:27"
$ws.Range("B147").Value = "11:08"
$ws.Range("C147").Value = "225_C ROCA-H SUR"
$ws.Range("D147").Value = 106
$ws.Range("E147").Value = "LP1912"
$ws.Range("A148").Value = "09:22:27"
$ws.Range("B148").Value = "11:09"
$ws.Range("C148").Value = "14_ABASTO"
$ws.Range("D148").Value = 107
$ws.Range("E148").Value = "LP1912"
$ws.Range("A149").Value = "09:22:27"
$ws.Range("B149").Value = "11:09"
$ws.Range("C149").Value = "17_ROMERO"
$ws.Range("D149").Value = 107
$ws.Range("E149").Value = "LP1912"
$ws.Range("A150").Value = "09:22:27"
$ws.Range("B150").Value = "11:19"
$ws.Range("C150").Value = "215C_EL PATO"
$ws.Range("D150").Value = 117
$ws.Range("E150").Value = "LP1912"
$ws.Range("A151").Value = "09:22:27"
$ws.Range("B151").Value = "11:20"
$ws.Range("C151").Value = "11_ETCHEVERRY"
$ws.Range("D151").Value = 118
$ws.Range("E151").Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 09:22:27"
$ws.Range("A3").Value = "Total filas: 23"
$ws.Range("A27").Value = "09:22:27"
$ws.Range("D27").Value = 42
$ws.Range("A28").Value = "09:22:27"
$ws.Range("B28").Value = "11:19"
$ws.Range("C28").Value = "215C_EL PATO"
$ws.Range("D28").Value = 117
$ws.Range("E28").Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 09:22:27"
$ws.Range("A17").Value = "09:22:27"
$ws.Range("D17").Value = 34
$ws.Range("A19").Value = "09:22:27"
$ws.Range("D19").Value = 49
$ws.Range("A21").Value = "09:22:27"
$ws.Range("D21").Value = 60
